# Auto-generated edit script: applies the "Sheets via scheduled runner" data refresh
# to the Lamia_Profits workbook (per-sheet leve-profit recompute).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 3
$ws.Range("H3").Value = 0
$ws.Range("J3").Value = 0
$ws.Range("L3").Value = 0
$ws.Range("N3").ClearContents()

# Row 15
$ws.Range("H15").Value = 1510.4375
$ws.Range("I15").Value = 1510.4375
$ws.Range("K15").Value = 4531.3125
$ws.Range("M15").Value = -4362.3125

# Row 102
$ws.Range("H102").Value = 0
$ws.Range("J102").Value = 0
$ws.Range("L102").Value = 0
$ws.Range("N102").ClearContents()

# Row 106
$ws.Range("H106").Value = 9697.333000000001
$ws.Range("I106").Value = 3990
$ws.Range("K106").Value = 3990
$ws.Range("M106").Value = -3359

$ws = $wb.Worksheets.Item("ARM")
# Row 61
$ws.Range("H61").Value = 3675.1628
$ws.Range("I61").Value = 3667.4285
$ws.Range("J61").Value = 4000
$ws.Range("K61").Value = 3667.4285
$ws.Range("L61").Value = 4000
$ws.Range("M61").Value = -3455.4285
$ws.Range("N61").Value = -4424

# Row 132
$ws.Range("H132").Value = 3428.8647
$ws.Range("I132").Value = 2325.182
$ws.Range("J132").Value = 12534.25
$ws.Range("K132").Value = 6975.545999999999
$ws.Range("L132").Value = 37602.75
$ws.Range("M132").Value = -4445.545999999999
$ws.Range("N132").Value = -42662.75

# Row 136
$ws.Range("H136").Value = 3675.1628
$ws.Range("I136").Value = 3667.4285
$ws.Range("J136").Value = 4000
$ws.Range("K136").Value = 11002.2855
$ws.Range("L136").Value = 12000
$ws.Range("M136").Value = -8452.2855
$ws.Range("N136").Value = -17100

$ws = $wb.Worksheets.Item("BSM")
# Row 57
$ws.Range("H57").Value = 116137.5
$ws.Range("J57").Value = 116137.5
$ws.Range("L57").Value = 116137.5
$ws.Range("N57").Value = -117577.5

# Row 134
$ws.Range("H134").Value = 1853.4584
$ws.Range("I134").Value = 1281.2609
$ws.Range("J134").Value = 15014
$ws.Range("K134").Value = 3843.7827
$ws.Range("L134").Value = 45042
$ws.Range("M134").Value = -1308.7827
$ws.Range("N134").Value = -50112

# Row 136
$ws.Range("H136").Value = 116137.5
$ws.Range("J136").Value = 116137.5
$ws.Range("L136").Value = 116137.5
$ws.Range("N136").Value = -126337.5

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 36192.79
$ws.Range("I31").Value = 4583.75
$ws.Range("K31").Value = 4583.75
$ws.Range("M31").Value = -4288.75

# Row 34
$ws.Range("H34").Value = 36192.79
$ws.Range("I34").Value = 4583.75
$ws.Range("K34").Value = 4583.75
$ws.Range("M34").Value = -4381.75

# Row 107
$ws.Range("H107").Value = 1906.7858
$ws.Range("I107").Value = 2149
$ws.Range("J107").Value = 1725.125
$ws.Range("K107").Value = 2149
$ws.Range("L107").Value = 1725.125
$ws.Range("M107").Value = -229
$ws.Range("N107").Value = -5565.125

# Row 132
$ws.Range("H132").Value = 2167.6206
$ws.Range("I132").Value = 1437.375
$ws.Range("K132").Value = 4312.125
$ws.Range("M132").Value = -1782.125

$ws = $wb.Worksheets.Item("CUL")
# Row 75
$ws.Range("H75").Value = 62504480
$ws.Range("I75").Value = 142857490
$ws.Range("K75").Value = 428572470
$ws.Range("M75").Value = -428571472

# Row 78
$ws.Range("H78").Value = 62504480
$ws.Range("I78").Value = 142857490
$ws.Range("K78").Value = 1285717410
$ws.Range("M78").Value = -1285712418

# Row 113
$ws.Range("H113").Value = 1113.75
$ws.Range("I113").Value = 932.8
$ws.Range("J113").Value = 1243
$ws.Range("K113").Value = 2798.4
$ws.Range("L113").Value = 3729
$ws.Range("M113").Value = -628.3999999999996
$ws.Range("N113").Value = -8069

$ws = $wb.Worksheets.Item("GSM")
# Row 43
$ws.Range("H43").Value = 2304.5
$ws.Range("I43").Value = 2304.5
$ws.Range("K43").Value = 2304.5
$ws.Range("M43").Value = -2153.5

# Row 113
$ws.Range("H113").Value = 5913.7896
$ws.Range("I113").Value = 5764.5835
$ws.Range("J113").Value = 6169.5713
$ws.Range("K113").Value = 5764.5835
$ws.Range("L113").Value = 6169.5713
$ws.Range("M113").Value = -3594.5835
$ws.Range("N113").Value = -10509.5713

# Row 126
$ws.Range("H126").Value = 3634.9333
$ws.Range("J126").Value = 6528.25
$ws.Range("L126").Value = 19584.75
$ws.Range("N126").Value = -24524.75

# Row 132
$ws.Range("H132").Value = 5947.6665
$ws.Range("I132").Value = 3402.9333
$ws.Range("K132").Value = 10208.7999
$ws.Range("M132").Value = -7678.7999

$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 8730.944
$ws.Range("I7").Value = 4316.364
$ws.Range("J7").Value = 15668.143
$ws.Range("K7").Value = 4316.364
$ws.Range("L7").Value = 15668.143
$ws.Range("M7").Value = -4204.364
$ws.Range("N7").Value = -15892.143

# Row 22
$ws.Range("H22").Value = 7924.846
$ws.Range("I22").Value = 1494.6666
$ws.Range("J22").Value = 9853.9
$ws.Range("K22").Value = 1494.6666
$ws.Range("L22").Value = 9853.9
$ws.Range("M22").Value = -1199.6666
$ws.Range("N22").Value = -10443.9

# Row 27
$ws.Range("H27").Value = 7924.846
$ws.Range("I27").Value = 1494.6666
$ws.Range("J27").Value = 9853.9
$ws.Range("K27").Value = 1494.6666
$ws.Range("L27").Value = 9853.9
$ws.Range("M27").Value = -1387.6666
$ws.Range("N27").Value = -10067.9

# Row 42
$ws.Range("H42").Value = 20000
$ws.Range("I42").Value = 20000
$ws.Range("K42").Value = 20000
$ws.Range("M42").Value = -19437

# Row 46
$ws.Range("H46").Value = 2681.7856
$ws.Range("I46").Value = 747.5
$ws.Range("J46").Value = 3004.1667
$ws.Range("K46").Value = 747.5
$ws.Range("L46").Value = 3004.1667
$ws.Range("M46").Value = -559.5
$ws.Range("N46").Value = -3380.1667

# Row 49
$ws.Range("H49").Value = 20000
$ws.Range("I49").Value = 20000
$ws.Range("K49").Value = 20000
$ws.Range("M49").Value = -19853

# Row 55
$ws.Range("H55").Value = 2274692.8
$ws.Range("J55").Value = 3519.4285
$ws.Range("L55").Value = 3519.4285
$ws.Range("N55").Value = -3865.4285

# Row 126
$ws.Range("H126").Value = 8730.944
$ws.Range("I126").Value = 4316.364
$ws.Range("J126").Value = 15668.143
$ws.Range("K126").Value = 12949.092
$ws.Range("L126").Value = 47004.429
$ws.Range("M126").Value = -10479.092
$ws.Range("N126").Value = -51944.429

# Row 136
$ws.Range("H136").Value = 7972.1577
$ws.Range("I136").Value = 1902.1111
$ws.Range("K136").Value = 5706.3333
$ws.Range("M136").Value = -3156.3333

$ws = $wb.Worksheets.Item("WVR")
# Row 132
$ws.Range("H132").Value = 5626.447
$ws.Range("I132").Value = 3357.7837
$ws.Range("J132").Value = 14020.5
$ws.Range("K132").Value = 10073.3511
$ws.Range("L132").Value = 42061.5
$ws.Range("M132").Value = -7543.3511
$ws.Range("N132").Value = -47121.5

# Row 136
$ws.Range("H136").Value = 2090.7856
$ws.Range("I136").Value = 1493.7307
$ws.Range("J136").Value = 9852.5
$ws.Range("K136").Value = 4481.1921
$ws.Range("L136").Value = 29557.5
$ws.Range("M136").Value = -1931.1921
